$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Add a new column "Tipo de Apoio" with support-type values for each node
$ws.Range("F1").Value = "Tipo de Apoio"
$ws.Range("F2").Value = 2
$ws.Range("F3").Value = 1
$ws.Range("F4").Value = 0

# Match the formatting of the existing numeric columns (A:C) - centered values
$ws.Range("F2:F4").HorizontalAlignment = -4108

# Update the active selection to the newly added cell, as happens after data entry
$ws.Range("F4").Select()
